# Update the "Parámetros" worksheet parameters:
#  - Número de trabajadores Full-Time (C3): 20 -> 15
#  - Número de trabajadores Part-Time (C4): 20 -> 10
#  - Turno (A o B) (C5): A -> B
# and move the cell cursor/selection to D4 (was L22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = "B"

$ws.Range("D4").Select()
